$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 21; $row++) {
    $cell = $ws.Cells.Item($row, 17)  # Column Q is the 17th column
    $cell.Value = $cell.Value2 * 2
}
